$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# "Ready for handoff" rows (4-7) in both the zh-cn and de-de sheets were
# re-generated for handoff: Priority moves from "low" to "ht", and the
# Latest Handoff Datetime is bumped to the new generation timestamp.
foreach ($row in 4..7) {
    $wsZh.Range("E$row").Value = "ht"
    $wsZh.Range("H$row").Value = "2016-08-16 22:29:37"

    $wsDe.Range("E$row").Value = "ht"
    $wsDe.Range("H$row").Value = "2016-08-16 22:29:43"

    # Overview's "Latest HO Xliff Generate Date" mirrors the newest
    # per-language handoff datetime (de-de is the later one here).
    $wsOverview.Range("G$row").Value = "2016-08-16 22:29:43"
}
